$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-01-27 Saturday" "2024-01-28 Sunday"

Replace-Text "46×58=" "80×59="
Replace-Text "74×34=" "14×47="
Replace-Text "29×75=" "48×17="
Replace-Text "77×37=" "81×53="
Replace-Text "11×78=" "50×11="
Replace-Text "76×41=" "53×88="
Replace-Text "66×54=" "54×97="
Replace-Text "48×37=" "48×90="
Replace-Text "92×97=" "53×84="
Replace-Text "25×77=" "22×79="
Replace-Text "65×29=" "96×81="
Replace-Text "82×16=" "85×92="
Replace-Text "64×49=" "20×17="
Replace-Text "96×35=" "86×96="
Replace-Text "11×46=" "60×28="
Replace-Text "88×97=" "79×91="
Replace-Text "96×18=" "94×91="
Replace-Text "59×82=" "96×83="
Replace-Text "97×42=" "13×19="
Replace-Text "53×32=" "36×53="
Replace-Text "69×70=" "98×13="
Replace-Text "26×82=" "72×30="
Replace-Text "94×46=" "79×37="
Replace-Text "50×91=" "48×98="
Replace-Text "95×30=" "41×35="
